# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - index 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2006
$ws1.Range("F7").Value  = 6275
$ws1.Range("F9").Value  = 1872
$ws1.Range("F10").Value = 496
$ws1.Range("F16").Value = 7386
$ws1.Range("F23").Value = 1
$ws1.Range("F26").Value = 45
$ws1.Range("F28").Value = 1661
$ws1.Range("F29").Value = 782
$ws1.Range("F30").Value = 344
$ws1.Range("F33").Value = 69

# Sheet "本地生活" (Local life) - index 3
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 256

# Sheet "全部类型" (All types) - index 4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 2006
$ws4.Range("F11").Value = 256
$ws4.Range("F13").Value = 6275
$ws4.Range("F15").Value = 1872
$ws4.Range("F18").Value = 496
$ws4.Range("F24").Value = 7386
$ws4.Range("F32").Value = 45
$ws4.Range("F34").Value = 1661
$ws4.Range("F36").Value = 344
